# KIBON-238: Tagi Code in einem Commit löschen
#
# The "Tagi" offering was dropped from the app, so its column in the
# "Benutzer" report (header "Tagi" in row 5 / placeholder "{isTagi}" in
# row 6) is no longer needed. Select the whole column and delete it
# outright, shifting every column to its right one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("L").Select() | Out-Null
$ws.Columns("L").Delete() | Out-Null
